# Quarterly income-statement roll-forward:
#  - drop the oldest quarter (column D) and shift every quarter one column
#    to the left (D<-E, E<-F, ... L<-M)
#  - populate the freed-up last column (M) with the new quarter
#  - a handful of historical figures were recomputed by the updated
#    read_price algorithm (same column position before the shift, i.e.
#    the old column J -> new column I for rows 19/20/22/24/25/27)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$dataCols = @("D","E","F","G","H","I","J","K","L","M")

function Shift-RowLeft($row) {
    for ($i = 0; $i -lt $dataCols.Length - 1; $i++) {
        $src = $dataCols[$i + 1]
        $dst = $dataCols[$i]
        $ws.Range($dst + $row).Value = $ws.Range($src + $row).Value()
    }
}

# ---- Row 8: quarter-ending labels ----
Shift-RowLeft 8
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# ---- Row 9: publish dates ----
Shift-RowLeft 9
$ws.Range("I9").Value = "1402-02-29 (8)"
$ws.Range("M9").Value = "1402-02-29"

# ---- Row 11: فروش (Sales) ----
Shift-RowLeft 11
$ws.Range("M11").Value = 863542430

# ---- Row 12: بهای تمام شده کالای فروش رفته (COGS) ----
Shift-RowLeft 12
$ws.Range("M12").Value = -743938052

# ---- Row 13: سود (زیان) ناخالص (Gross profit) ----
Shift-RowLeft 13
$ws.Range("M13").Value = 119604378

# ---- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ----
Shift-RowLeft 14
$ws.Range("M14").Value = -9270177

# ---- Row 15: هزینه کاهش ارزش دریافتنی‌ها (Impairment expense) - all zero ----
Shift-RowLeft 15
$ws.Range("M15").Value = 0

# ---- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense) ----
Shift-RowLeft 16
$ws.Range("M16").Value = -160675

# ---- Row 17: سود (زیان) عملیاتی (Operating profit) ----
Shift-RowLeft 17
$ws.Range("M17").Value = 110173526

# ---- Row 18: هزینه های مالی (Finance costs) ----
Shift-RowLeft 18
$ws.Range("M18").Value = -4482865

# ---- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating items) ----
# recalculated figure lands in the same slot as before the shift (old J -> new I)
Shift-RowLeft 19
$ws.Range("I19").Value = -1162243
$ws.Range("M19").Value = 15788063

# ---- Row 20: سود خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ----
Shift-RowLeft 20
$ws.Range("I20").Value = 76942704
$ws.Range("M20").Value = 121478724

# ---- Row 21: مالیات (Tax) ----
Shift-RowLeft 21
$ws.Range("M21").Value = -1513686

# ---- Row 22: سود خالص عملیات در حال تداوم (Net profit from continuing ops) ----
Shift-RowLeft 22
$ws.Range("I22").Value = 69571420
$ws.Range("M22").Value = 119965038

# ---- Row 23: سود عملیات متوقف شده (Discontinued ops) - all zero ----
Shift-RowLeft 23
$ws.Range("M23").Value = 0

# ---- Row 24: سود (زیان) خالص (Net profit) ----
Shift-RowLeft 24
$ws.Range("I24").Value = 69571420
$ws.Range("M24").Value = 119965038

# ---- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ----
Shift-RowLeft 25
$ws.Range("I25").Value = 366
$ws.Range("M25").Value = 528

# ---- Row 26: سرمایه (Capital) ----
Shift-RowLeft 26
$ws.Range("M26").Value = 227000000

# ---- Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS on latest capital) ----
# fully recalculated against the new (latest) capital value, not a simple shift
$ws.Range("D27").Value = 188
$ws.Range("E27").Value = 275
$ws.Range("F27").Value = 113
$ws.Range("G27").Value = 86
$ws.Range("H27").Value = 224
$ws.Range("I27").Value = 306
$ws.Range("J27").Value = 578
$ws.Range("K27").Value = 352
$ws.Range("L27").Value = 223
$ws.Range("M27").Value = 528
